$d = $word.ActiveDocument

# Locate the paragraph ending in ". Simile con gente" and insert a brand
# new paragraph right after it with the new bug report line.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Simile con gente") {
        $target = $p
        break
    }
}

$target.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($i + 1)

# Append a sentinel trailing character so the bookmark we add below is not
# collapsed exactly on the paragraph-end boundary (insert position), then
# strip the sentinel back out once the bookmark is anchored.
$newPara.Range.Text = "-Cuando solicitas gente hace un request a amigos!!! Que obviamente falla porque no estas loguineadoX"

$anchorPos = $newPara.Range.End - 2
$anchorRange = $d.Range($anchorPos, $anchorPos)
$d.Bookmarks.Add("_GoBack", $anchorRange)

$sentinelRange = $d.Range($newPara.Range.End - 2, $newPara.Range.End - 1)
$sentinelRange.Delete()
